# Auto-generated edit script: updates market price data cells across all 8 crafting-job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect a scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1237.375
$ws.Range("I28").Value = 1271.2858
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 1271.2858
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = -786.2858000000001
$ws.Range("N28").Value = -1970

$ws.Range("H43").Value = 3477.111
$ws.Range("J43").Value = 3374.5
$ws.Range("L43").Value = 3374.5
$ws.Range("N43").Value = -3512.5

$ws.Range("H62").Value = 6793.5415
$ws.Range("I62").Value = 6684.6113
$ws.Range("K62").Value = 6684.6113
$ws.Range("M62").Value = -6060.6113

$ws.Range("H65").Value = 6793.5415
$ws.Range("I65").Value = 6684.6113
$ws.Range("K65").Value = 33423.0565
$ws.Range("M65").Value = -30303.0565

$ws.Range("H103").Value = 556.7895
$ws.Range("J103").Value = 789.1818
$ws.Range("L103").Value = 2367.5454
$ws.Range("N103").Value = -3539.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1161.7576
$ws.Range("I2").Value = 1188.25
$ws.Range("K2").Value = 1188.25
$ws.Range("M2").Value = -1075.25

$ws.Range("H5").Value = 4816.909
$ws.Range("I5").Value = 8606.5
$ws.Range("J5").Value = 269.4
$ws.Range("K5").Value = 8606.5
$ws.Range("L5").Value = 269.4
$ws.Range("M5").Value = -8494.5
$ws.Range("N5").Value = -493.4

$ws.Range("H32").Value = 11253.5
$ws.Range("I32").Value = 12147.929
$ws.Range("K32").Value = 12147.929
$ws.Range("M32").Value = -11860.929

$ws.Range("H74").Value = 3144.7837
$ws.Range("J74").Value = 6348.5
$ws.Range("L74").Value = 6348.5
$ws.Range("N74").Value = -8096.5

$ws.Range("H77").Value = 3144.7837
$ws.Range("J77").Value = 6348.5
$ws.Range("L77").Value = 31742.5
$ws.Range("N77").Value = -40478.5

$ws.Range("H116").Value = 1161.7576
$ws.Range("I116").Value = 1188.25
$ws.Range("K116").Value = 1188.25
$ws.Range("M116").Value = 1105.75

$ws.Range("H122").Value = 3057.2
$ws.Range("I122").Value = 2640.9583
$ws.Range("J122").Value = 4722.1665
$ws.Range("K122").Value = 7922.874899999999
$ws.Range("L122").Value = 14166.4995
$ws.Range("M122").Value = -5472.874899999999
$ws.Range("N122").Value = -19066.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1161.7576
$ws.Range("I3").Value = 1188.25
$ws.Range("K3").Value = 1188.25
$ws.Range("M3").Value = -1074.25

$ws.Range("H4").Value = 4816.909
$ws.Range("I4").Value = 8606.5
$ws.Range("J4").Value = 269.4
$ws.Range("K4").Value = 8606.5
$ws.Range("L4").Value = 269.4
$ws.Range("M4").Value = -8491.5
$ws.Range("N4").Value = -499.4

$ws.Range("H20").Value = 2766.8262
$ws.Range("J20").Value = 3017.3
$ws.Range("L20").Value = 3017.3
$ws.Range("N20").Value = -3511.3

$ws.Range("H22").Value = 781.6667
$ws.Range("I22").Value = 781.6667
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 781.6667
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -608.6667
$ws.Range("N22").Value = ""

$ws.Range("H105").Value = 7050
$ws.Range("I105").Value = 4100
$ws.Range("K105").Value = 4100
$ws.Range("M105").Value = -2353

$ws.Range("H107").Value = 2500
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = ""

$ws.Range("H134").Value = 7054.7915
$ws.Range("I134").Value = 2507.7646
$ws.Range("K134").Value = 7523.293799999999
$ws.Range("M134").Value = -4988.293799999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2271.4167
$ws.Range("J58").Value = 1998.2
$ws.Range("L58").Value = 1998.2
$ws.Range("N58").Value = -2404.2

$ws.Range("H134").Value = 2923.2856
$ws.Range("I134").Value = 1292.6444
$ws.Range("K134").Value = 3877.933199999999
$ws.Range("M134").Value = -1342.933199999999

$ws.Range("H136").Value = 2271.4167
$ws.Range("J136").Value = 1998.2
$ws.Range("L136").Value = 5994.6
$ws.Range("N136").Value = -11094.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 326.16666
$ws.Range("I38").Value = 331.77777
$ws.Range("J38").Value = 309.33334
$ws.Range("K38").Value = 995.33331
$ws.Range("L38").Value = 928.0000200000001
$ws.Range("M38").Value = -648.33331
$ws.Range("N38").Value = -1622.00002

$ws.Range("H70").Value = 312
$ws.Range("I70").Value = 312
$ws.Range("K70").Value = 936
$ws.Range("M70").Value = -621

$ws.Range("H73").Value = 312
$ws.Range("I73").Value = 312
$ws.Range("K73").Value = 936
$ws.Range("M73").Value = 156

$ws.Range("H75").Value = 2187.6667
$ws.Range("I75").Value = 1076.5714
$ws.Range("J75").Value = 3743.2
$ws.Range("K75").Value = 3229.7142
$ws.Range("L75").Value = 11229.6
$ws.Range("M75").Value = -2231.7142
$ws.Range("N75").Value = -13225.6

$ws.Range("H78").Value = 2187.6667
$ws.Range("I78").Value = 1076.5714
$ws.Range("J78").Value = 3743.2
$ws.Range("K78").Value = 9689.142600000001
$ws.Range("L78").Value = 33688.8
$ws.Range("M78").Value = -4697.142600000001
$ws.Range("N78").Value = -43672.8

$ws.Range("H97").Value = 1939.5
$ws.Range("J97").Value = 786.7143
$ws.Range("L97").Value = 2360.1429
$ws.Range("N97").Value = -3352.1429

$ws.Range("H103").Value = 568
$ws.Range("J103").Value = 621.8333
$ws.Range("L103").Value = 1865.4999
$ws.Range("N103").Value = -3623.4999

$ws.Range("H107").Value = 4711
$ws.Range("J107").Value = 5573.9
$ws.Range("L107").Value = 16721.7
$ws.Range("N107").Value = -20561.7

$ws.Range("H109").Value = 2151.8333
$ws.Range("I109").Value = 1282.7
$ws.Range("K109").Value = 3848.1
$ws.Range("M109").Value = -2808.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8009.4
$ws.Range("I70").Value = 7761.75
$ws.Range("K70").Value = 7761.75
$ws.Range("M70").Value = -7491.75

$ws.Range("H73").Value = 8009.4
$ws.Range("I73").Value = 7761.75
$ws.Range("K73").Value = 7761.75
$ws.Range("M73").Value = -6825.75

$ws.Range("H103").Value = 38000
$ws.Range("J103").Value = 38000
$ws.Range("L103").Value = 38000
$ws.Range("N103").Value = -40344

$ws.Range("H122").Value = 1851.6666
$ws.Range("I122").Value = 1447.05
$ws.Range("K122").Value = 4341.15
$ws.Range("M122").Value = -1891.15

$ws.Range("H132").Value = 6805288.5
$ws.Range("I132").Value = 8132332
$ws.Range("K132").Value = 24396996
$ws.Range("M132").Value = -24394466

$ws.Range("H133").Value = 70776.664
$ws.Range("J133").Value = 70776.664
$ws.Range("L133").Value = 70776.664
$ws.Range("N133").Value = -80896.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7328.3076
$ws.Range("I7").Value = 5644.5
$ws.Range("J7").Value = 8771.571
$ws.Range("K7").Value = 5644.5
$ws.Range("L7").Value = 8771.571
$ws.Range("M7").Value = -5532.5
$ws.Range("N7").Value = -8995.571

$ws.Range("H93").Value = 3078.8
$ws.Range("I93").Value = 2465
$ws.Range("J93").Value = 3999.5
$ws.Range("K93").Value = 2465
$ws.Range("L93").Value = 3999.5
$ws.Range("M93").Value = -1217
$ws.Range("N93").Value = -6495.5

$ws.Range("H112").Value = 88499.5
$ws.Range("J112").Value = 88499.5
$ws.Range("L112").Value = 88499.5
$ws.Range("N112").Value = -91453.5

$ws.Range("H122").Value = 3432.0625
$ws.Range("I122").Value = 2455
$ws.Range("J122").Value = 7666
$ws.Range("K122").Value = 7365
$ws.Range("L122").Value = 22998
$ws.Range("M122").Value = -4915
$ws.Range("N122").Value = -27898

$ws.Range("H126").Value = 7328.3076
$ws.Range("I126").Value = 5644.5
$ws.Range("J126").Value = 8771.571
$ws.Range("K126").Value = 16933.5
$ws.Range("L126").Value = 26314.713
$ws.Range("M126").Value = -14463.5
$ws.Range("N126").Value = -31254.713

$ws.Range("H128").Value = 62498
$ws.Range("J128").Value = 62498
$ws.Range("L128").Value = 62498
$ws.Range("N128").Value = -72458

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 89998.664
$ws.Range("J46").Value = 89998.664
$ws.Range("L46").Value = 89998.664
$ws.Range("N46").Value = -90460.664

$ws.Range("H94").Value = 32020.2
$ws.Range("J94").Value = 32020.2
$ws.Range("L94").Value = 32020.2
$ws.Range("N94").Value = -33822.2

$ws.Range("H134").Value = 89998.664
$ws.Range("J134").Value = 89998.664
$ws.Range("L134").Value = 269995.992
$ws.Range("N134").Value = -275065.992

$ws.Range("H136").Value = 30035.143
$ws.Range("I136").Value = 34541
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 103623
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -101073
$ws.Range("N136").Value = -14100
